$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46; this shifts the existing rows 46-111 down to 47-112
$ws.Rows.Item(46).Insert()

# The newly inserted row 46 is blank. Populate it by copying the row that is now
# row 47 (the original row 46 data), so all the "static" columns (A,B,C,E,F,G,H,I,N,O,Q,R)
# carry over identically, matching the new entry's metadata.
$ws.Range("A47:R47").Copy() | Out-Null
$ws.Range("A46").PasteSpecial() | Out-Null

# Now overwrite the columns that hold the new weekly record's own data.
$ws.Range("D46").Value = 44799
$ws.Range("J46").Value = 50
$ws.Range("K46").Value = 11000
$ws.Range("L46").Value = 12000
$ws.Range("M46").Value = 11600
$ws.Range("P46").Value = 193
